$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 226-241: team name slugs + refreshed odds ---
$updates = @(
    @{ Row = 226; B = "tennessee-titans";       C = "san-francisco-49ers";   D = 150;  E = -170; F = 3.5;  G = 45 },
    @{ Row = 227; B = "green-bay-packers";      C = "cleveland-browns";      D = -365; E = 280;  F = 7.5;  G = 47 },
    @{ Row = 228; B = "arizona-cardinals";      C = "indianapolis-colts";    D = -165; E = 145;  F = 3;    G = 48 },
    @{ Row = 229; B = "atlanta-falcons";        C = "detroit-lions";         D = -365; E = 280;  F = 7.5;  G = 43 },
    @{ Row = 230; B = "cincinnati-bengals";     C = "baltimore-ravens";      D = -380; E = 290;  F = 7.5;  G = 43 },
    @{ Row = 231; B = "minnesota-vikings";      C = "los-angeles-rams";      D = 150;  E = -170; F = 3;    G = 49 },
    @{ Row = 232; B = "new-england-patriots";   C = "buffalo-bills";         D = -120; E = 100;  F = 1;    G = 43.5 },
    @{ Row = 233; B = "new-york-jets";          C = "jacksonville-jaguars";  D = -130; E = 110;  F = 2;    G = 43 },
    @{ Row = 234; B = "philadelphia-eagles";    C = "new-york-giants";       D = -550; E = 400;  F = 11;   G = 41 },
    @{ Row = 235; B = "carolina-panthers";      C = "tampa-bay-buccaneers";  D = 400;  E = -550; F = 11.5; G = 44 },
    @{ Row = 236; B = "houston-texans";         C = "los-angeles-chargers";  D = 460;  E = -650; F = 12.5; G = 45.5 },
    @{ Row = 237; B = "seattle-seahawks";       C = "chicago-bears";         D = -335; E = 260;  F = 7;    G = 41.5 },
    @{ Row = 238; B = "kansas-city-chiefs";     C = "pittsburgh-steelers";   D = -450; E = 340;  F = 10;   G = 44 },
    @{ Row = 239; B = "las-vegas-raiders";      C = "denver-broncos";        D = -105; E = -115; F = 1;    G = 41.5 },
    @{ Row = 240; B = "dallas-cowboys";         C = "washington-football-team"; D = -475; E = 350; F = 10; G = 46 },
    @{ Row = 241; B = "new-orleans-saints";     C = "miami-dolphins";        D = 135;  E = -155; F = 3;    G = 37 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("B$r").Value = $u.B
    $ws.Range("C$r").Value = $u.C
    $ws.Range("D$r").Value = $u.D
    $ws.Range("E$r").Value = $u.E
    $ws.Range("F$r").Value = $u.F
    $ws.Range("G$r").Value = $u.G
}

# --- Remove the two trailing duplicate rows (242 and 243) ---
$ws.Rows("242:243").Delete() | Out-Null
